# Updated cryptos list with GitHub Actions — refresh Price (D) / Volume(1h) (E)
# columns for the coin rows on Sheet1. Values that read as plain numbers
# (e.g. "1.013") are written with a leading apostrophe so Excel keeps them
# as literal text instead of coercing them to a Double, then the cell
# style is put back to "Normal" so no stray number-format is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.964.22'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '1.843.75'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("D4").Value = "'1.013"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").Value = "'1.011"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").Value = "'308.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").Value = "'0.4772"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.91%  '
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("D9").Value = "'0.07212"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").Value = "'0.9291"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("D11").Value = "'19.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").Value = '1.858.62'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").Value = "'5.405"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").Value = "'6.451"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.76%  '
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").Value = "'0.000008654"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("D19").Value = "'1.012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = '27.006.96'
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("D22").Value = "'5.067"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").Value = "'10.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").Value = "'1.947"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").Value = "'152.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("D27").Value = "'2.013"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.88%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").Value = "'4.966"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("D30").Value = "'0.08861"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").Value = "'3.313"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.72%  '
$ws.Range("D32").Value = "'1.180"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("D33").Value = "'0.7416"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").Value = "'4.494"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").Value = "'2.701"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.67%  '
$ws.Range("D36").Value = "'1.112"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.86%  '
$ws.Range("E37").Value = '  +0.95%  '
$ws.Range("D38").Value = "'0.05251"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'2.965"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("D40").Value = "'0.5256"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").Value = "'7.011"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D43").Value = "'8.283"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("D45").Value = "'0.4732"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").Value = "'1.012"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").Value = "'101.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").Value = "'65.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("D50").Value = "'0.06073"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("D51").Value = "'0.8887"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.19%  '
